$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G column values
$ws.Range("G5").Value = 54.6
$ws.Range("G8").Value = 51
$ws.Range("G9").Value = 51
$ws.Range("G10").Value = 51
$ws.Range("G11").Value = 51
$ws.Range("G12").Value = 51
$ws.Range("G13").Value = 51
$ws.Range("G14").Value = 51
$ws.Range("G15").Value = 51
$ws.Range("G16").Value = 51
$ws.Range("G17").Value = 51
$ws.Range("G18").Value = 51
$ws.Range("G24").Value = 51
$ws.Range("G25").Value = 51
$ws.Range("G26").Value = 51

# Update the active cell selection to H25
$ws.Range("H25").Select()
